$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2EXT04_DNA")

# Fill in example values for row 2
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "Total RNA"
$ws.Range("C2").Value = "NCIT"
$ws.Range("D2").Value = "https://bioregistry.io/NCIT:C163995"
$ws.Range("E2").Value = "200"
$ws.Range("F2").Value = "milligram"
$ws.Range("G2").Value = "UO"
$ws.Range("H2").Value = "https://bioregistry.io/UO:0000022"
$ws.Range("I2").Value = "QIAGEN RNEasy"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "QIAGEN RNEasy Buffer 2"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "200"
$ws.Range("P2").Value = "microliter"
$ws.Range("Q2").Value = "UO"
$ws.Range("R2").Value = "https://bioregistry.io/UO:0000101"
$ws.Range("S2").Value = ""

# Delete rows 3-5 (remove extra example rows)
$ws.Range("A3:S5").Delete()
